$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Overview" - rows for 7f919e97 (row 4) and 8022417e (row 5) move from
# "Ready for handoff" to "Handed back: in sync with en-US" for both the
# zh-cn (E) and de-de (F) status columns.
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E4").Value = "Handed back: in sync with en-US"
$wsOverview.Range("F4").Value = "Handed back: in sync with en-US"
$wsOverview.Range("E5").Value = "Handed back: in sync with en-US"
$wsOverview.Range("F5").Value = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------------
# Sheet "zh-cn" - mark 7f919e97 (row 4) and 8022417e (row 5) as handed back:
#   - Status (C) becomes "Handed back: in sync with en-US"
#   - Latest Target File (I) gets the source file name + hyperlink
#   - Latest Handback File (J) gets the handback xlf file name
#   - Latest Handback DateTime (K) gets the handback timestamp
# ---------------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Range("C4").Value = "Handed back: in sync with en-US"
$wsZhCn.Range("I4").Value = "7f919e97-f145-4d66-aedf-7bf6cb69e2f4.md"
$wsZhCn.Range("J4").Value = "7f919e97-f145-4d66-aedf-7bf6cb69e2f4.edfe68b67b637a248421e29a868a494329015ee9.zh-cn.xlf"
$wsZhCn.Range("K4").Value = "2016-10-19 12:26:21"
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I4"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/96d3a9616f5710b534273ffb9ee858e91bf85d23/e2e/7f919e97-f145-4d66-aedf-7bf6cb69e2f4.md", [Type]::Missing, [Type]::Missing, "7f919e97-f145-4d66-aedf-7bf6cb69e2f4.md")

$wsZhCn.Range("C5").Value = "Handed back: in sync with en-US"
$wsZhCn.Range("I5").Value = "8022417e-44ba-40ff-a407-d2fcf7c88ae4.md"
$wsZhCn.Range("J5").Value = "8022417e-44ba-40ff-a407-d2fcf7c88ae4.731d0ab17cfc6d0bf85f99ce1626667af748129e.zh-cn.xlf"
$wsZhCn.Range("K5").Value = "2016-10-19 12:26:21"
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I5"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/96d3a9616f5710b534273ffb9ee858e91bf85d23/e2e/8022417e-44ba-40ff-a407-d2fcf7c88ae4.md", [Type]::Missing, [Type]::Missing, "8022417e-44ba-40ff-a407-d2fcf7c88ae4.md")

# ---------------------------------------------------------------------------
# Sheet "de-de" - same handback bookkeeping as zh-cn, different timestamp and
# repo used for the hyperlink targets.
# ---------------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Range("C4").Value = "Handed back: in sync with en-US"
$wsDeDe.Range("I4").Value = "7f919e97-f145-4d66-aedf-7bf6cb69e2f4.md"
$wsDeDe.Range("J4").Value = "7f919e97-f145-4d66-aedf-7bf6cb69e2f4.edfe68b67b637a248421e29a868a494329015ee9.de-de.xlf"
$wsDeDe.Range("K4").Value = "2016-10-19 12:26:38"
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I4"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/3cee8b652e1d17d333c508687eff950b482571df/e2e/7f919e97-f145-4d66-aedf-7bf6cb69e2f4.md", [Type]::Missing, [Type]::Missing, "7f919e97-f145-4d66-aedf-7bf6cb69e2f4.md")

$wsDeDe.Range("C5").Value = "Handed back: in sync with en-US"
$wsDeDe.Range("I5").Value = "8022417e-44ba-40ff-a407-d2fcf7c88ae4.md"
$wsDeDe.Range("J5").Value = "8022417e-44ba-40ff-a407-d2fcf7c88ae4.731d0ab17cfc6d0bf85f99ce1626667af748129e.de-de.xlf"
$wsDeDe.Range("K5").Value = "2016-10-19 12:26:38"
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I5"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/3cee8b652e1d17d333c508687eff950b482571df/e2e/8022417e-44ba-40ff-a407-d2fcf7c88ae4.md", [Type]::Missing, [Type]::Missing, "8022417e-44ba-40ff-a407-d2fcf7c88ae4.md")
